$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "66.627.99"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.95%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.200.14"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.76%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "575.64"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.66%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "137.79"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  -15.68%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "3.198.51"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.01%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.516"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  -10.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  -14.39%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "6.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.38%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.470"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.57%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "35.55"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  -14.92%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000227"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.34%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.727.38"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "66.699.80"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.84%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "3.214.39"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.41%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.31%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "6.65"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  -14.18%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "488.93"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  -12.57%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "14.11"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  -13.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.707"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  -12.16%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  -15.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "80.91"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  -10.16%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "12.58"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.61%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.86%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 2)
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "27.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  -12.53%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 2)
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  -12.38%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "7.39"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  -9.38%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.14"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.77%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  -7.45%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.52%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 2)
$c.NumberFormat = "@"
$c.Value = "OKB"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "54.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.50%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 2)
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  -18.42%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 2)
$c.NumberFormat = "@"
$c.Value = "Bittensor"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "486.89"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -15.52%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "5.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -16.37%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.0410"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -8.77%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.0802"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.77%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  -13.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "8.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  -16.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.801.37"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  -9.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -15.82%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 2)
$c.NumberFormat = "@"
$c.Value = "USDe"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 2)
$c.NumberFormat = "@"
$c.Value = "TheGraph"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.247"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  -11.25%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "121.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  -6.29%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "24.56"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  -15.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = "Fetch.AI"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  -10.87%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0516"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -17.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.108"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -10.77%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "2.10"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  -20.89%  "
$c.Style = "Normal"
